$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista_dokumentów")

# Cells that only get the generic "PPP-NNN" placeholder keep the sheet's
# default (unstyled) look, so after writing the value we reset the style
# back to "Normalny" (this workbook's Normal cell style) to drop the
# centered/wrap formatting the B column otherwise inherits.
$plainCells = @(2, 4, 8, 14, 17, 18, 20, 21, 23, 24, 25, 26, 27, 28)
foreach ($r in $plainCells) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = "PPP-NNN"
    $cell.Style = "Normalny"
}

# Cells that received an actual KZ number keep the column's normal
# (centered) style, so a plain value assignment is enough.
$kzCells = @{
    3  = "5123, PPP-NNN"
    5  = "321, PPP-NNN"
    7  = "423, PPP-NNN"
    11 = "1523, PPP-NNN"
    13 = "1234, PPP-NNN"
    19 = "1235, PPP-NNN"
}
foreach ($r in $kzCells.Keys) {
    $ws.Cells.Item($r, 2).Value = $kzCells[$r]
}

# Column B widened slightly (raw OOXML width 9.140625 -> 10) + best-fit,
# matching Excel's "best fit" auto-resize after typing into the column.
$ws.Columns.Item(2).ColumnWidth = 9.140625
$ws.Columns.Item(2).BestFit = $true

# Selection moved from M13 to F11.
$ws.Range("F11").Select() | Out-Null
